# Daily attendance processing - 2026-02-01 15:39:14
#
# The "Recorded By" column (G) lists the user(s)/system accounts that
# recorded a session, separated by ", ". This pass normalizes the stored
# order of those names for a known set of previously-seen combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact "Recorded By" strings to normalize, mapped old -> new (order swapped).
$replacements = @{
    "system, backup@backdoor.com, System" = "backup@backdoor.com, system, System"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "admin@admin.com, System"             = "System, admin@admin.com"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G holds "Recorded By".
$col = 7

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Value2
    if ($null -eq $current) { continue }
    $text = [string]$current
    if ($replacements.ContainsKey($text)) {
        $cell.Value = $replacements[$text]
    }
}
